# My_readings.xlsx — "new variables, new issues, new resolves..."
# Fills in the previously-zeroed start_time (B) / finish_time (C) columns
# for a batch of reading sessions, and updates the sheet's scroll/selection
# state to where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (start_time, finish_time) as Excel time-of-day serial fractions
$times = @{
    2  = @(0.625,              0.65625)
    3  = @(0.708333333333333,  0.75)
    4  = @(0.791666666666667,  0.833333333333333)
    5  = @(0.840277777777778,  0.861111111111111)
    6  = @(0.75,               0.791666666666667)
    7  = @(0.791666666666667,  0.833333333333333)
    8  = @(0.583333333333333,  0.604166666666667)
    9  = @(0.666666666666667,  0.694444444444444)
    10 = @(0.75,               0.78125)
    14 = @(0.708333333333333,  0.729166666666667)
    15 = @(0.833333333333333,  0.854166666666667)
    16 = @(0.791666666666667,  0.833333333333333)
    17 = @(0.836805555555556,  0.857638888888889)
    18 = @(0.75,               0.763888888888889)
    19 = @(0.770833333333333,  0.8125)
    20 = @(0.5,                0.541666666666667)
    21 = @(0.583333333333333,  0.604166666666667)
    22 = @(0.708333333333333,  0.729166666666667)
    23 = @(0.833333333333333,  0.836805555555556)
    24 = @(0.840277777777778,  0.881944444444444)
    25 = @(0.708333333333333,  0.763888888888889)
    26 = @(0.791666666666667,  0.875)
    27 = @(0.833333333333333,  0.916666666666667)
}

foreach ($row in $times.Keys) {
    $pair = $times[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}

# Move the view/selection to where editing left off.
$ws.Range("A13").Select()
$ws.Range("F34").Select()
